$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.232.29'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '2.188.70'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''255.30'
$ws.Range('E5').Value = '  +4.14%  '
$ws.Range('D6').Value = '''0.629'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''68.32'
$ws.Range('E7').Value = '  -2.74%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = '''0.574'
$ws.Range('E9').Value = '  +2.87%  '
$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').Value = '''37.24'
$ws.Range('E10').Value = '  -4.98%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '''58.89'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('D13').Value = '''7.02'
$ws.Range('E13').Value = '  +3.75%  '
$ws.Range('D14').Value = '''0.104'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '2.514.22'
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').Value = '''0.874'
$ws.Range('E16').Value = '  +3.57%  '
$ws.Range('D17').Value = '''14.43'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').Value = '2.216.53'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').Value = '41.183.43'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').Value = '0.0₃0957'
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').Value = '''6.17'
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('D22').Value = '''72.13'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').Value = '''233.12'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('E24').Value = '  -3.48%  '
$ws.Range('D25').Value = '''11.81'
$ws.Range('E25').Value = '  +19.73%  '
$ws.Range('D26').Value = '''3.85'
$ws.Range('E26').Value = '  +5.10%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  +3.27%  '
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('D30').Value = '''169.18'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').Value = '''20.66'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('D33').Value = '''0.0749'
$ws.Range('E33').Value = '  +4.65%  '
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('D35').Value = '''5.46'
$ws.Range('E35').Value = '  +4.42%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '''4.16'
$ws.Range('E36').Value = '  +6.62%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').Value = '''25.94'
$ws.Range('E37').Value = '  +8.09%  '
$ws.Range('D38').Value = '''4.60'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '''0.0299'
$ws.Range('E39').Value = '  +8.20%  '
$ws.Range('D40').Value = '''2.20'
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('D41').Value = '''12.16'
$ws.Range('E41').Value = '  +12.84%  '
$ws.Range('D42').Value = '''5.66'
$ws.Range('E42').Value = '  -3.44%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').Value = '''63.35'
$ws.Range('E43').Value = '  -4.89%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').Value = '''4.93'
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('D45').Value = '''0.197'
$ws.Range('E45').Value = '  -4.97%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '''1.22'
$ws.Range('E46').Value = '  +11.66%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '''8.64'
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.101'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').Value = '''1.00'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').Value = '''1.18'
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('E51').Value = '  -9.41%  '
